$d = $word.ActiveDocument

# --- 1) Standardize the four job-title cells (first row, first column of
#        the first few experience tables) -------------------------------
# NOTE: Wrap=0 (wdFindStop) + Replace=1 (wdReplaceOne) keeps each Find
# strictly confined to the single cell Range it was scoped to; using
# wdReplaceAll here would leak into other tables that contain the same
# literal job-title text.

# Table 1: "Software Developer Team Lead" -> "Team Leader, Software Developer"
$cell1 = $d.Tables(1).Cell(1,1)
$cell1.Range.Find.Execute("Software Developer Team Lead", $true, $false, $false, $false, $false, $true, 0, $false, "Team Leader, Software Developer", 1)

# Table 2: "Senior Software Developer" -> "Senior Software Developer " (trailing space)
$cell2 = $d.Tables(2).Cell(1,1)
$cell2.Range.Find.Execute("Senior Software Developer", $true, $false, $false, $false, $false, $true, 0, $false, "Senior Software Developer ", 1)

# Table 3: "Software Developer Senior" -> "Senior Software Developer " (trailing space)
$cell3 = $d.Tables(3).Cell(1,1)
$cell3.Range.Find.Execute("Software Developer Senior", $true, $false, $false, $false, $false, $true, 0, $false, "Senior Software Developer ", 1)

# Table 5: "Software Developer Senior" -> "Senior Software Developer" (no trailing space)
$cell5 = $d.Tables(5).Cell(1,1)
$cell5.Range.Find.Execute("Software Developer Senior", $true, $false, $false, $false, $false, $true, 0, $false, "Senior Software Developer", 1)

# --- 2) Relocate the "_GoBack" bookmark from the name heading down to the
#        start of the first job-title cell (last edit point). -----------

$target = $d.Tables(1).Cell(1,1).Range
$insertionPoint = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $insertionPoint)
